$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header cell's style/format (AC1) onto the new header
# cells so that AD1:AF1 share the same bold/centered/bordered style as the
# rest of row 1, then set the header text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team record (Wins/Losses/Ties) for every data row (2-50).
for ($r = 2; $r -le 50; $r++) {
    $ws.Cells.Item($r, 30).Value = 81   # AD - Wins
    $ws.Cells.Item($r, 31).Value = 81   # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF - Ties
}

Write-Output "done"
